$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The target line of code (in the "swap" example) reads:
#     int* tmp = *a;
# and after the edit it should read:
#     int tmp =  a;
# i.e. the first "*" (right after "int", rendered with a lowered baseline
# position) is removed while keeping the following space, and the second
# "*" (a separate red-colored run right after " = ") is removed entirely.
# The now-empty spot where that second run used to live is where Word's
# "_GoBack" bookmark is relocated to (it previously sat at the very end of
# the document, after the closing "}" of the last code block).
# ---------------------------------------------------------------------------

# Locate the unique phrase containing both runs we need to touch.
$rng = $d.Content
$found = $rng.Find.Execute("int* tmp = *a;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target text 'int* tmp = *a;'"
}

$paraStart = $rng.Start

# --- Edit 1: "* " (before "tmp") becomes " " -------------------------------
# That run is exactly the 2 characters right after "int".
$starSpace = $d.Range($paraStart + 3, $paraStart + 5)
if ($starSpace.Text -ne "* ") {
    throw "Unexpected text at star+space range: [$($starSpace.Text)]"
}
$starSpace.Text = " "

# --- Edit 2: remove the red "*" run that sits right after " = " ------------
# Edit 1 shortened the paragraph by one character, so re-locate the line via
# Find rather than trusting the now-stale offsets computed before edit 1.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("tmp = *a;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not relocate 'tmp = *a;' line after first edit"
}

# Within the found phrase, locate the literal "*" character immediately
# preceding "a;" and remove it, then drop a _GoBack bookmark at that spot.
$phrase = $rng2.Text
$starIdx = $phrase.IndexOf("*")
if ($starIdx -lt 0) {
    throw "Could not find '*' inside phrase [$phrase]"
}
$starPos = $rng2.Start + $starIdx

# Remove any pre-existing "_GoBack" bookmark (Word only ever keeps a single
# one) and re-insert it right before the "*" we are about to delete. Doing
# this *before* deleting the run keeps the bookmark as a boundary between
# the deleted run and the following "a;" run, so that run is left untouched
# (no accidental merge / no spurious xml:space="preserve").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($starPos, $starPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$starRange = $d.Range($starPos, $starPos + 1)
if ($starRange.Text -ne "*") {
    throw "Unexpected text at star range: [$($starRange.Text)]"
}
$starRange.Text = ""
